$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.375.54"
Set-TextValue $ws.Range("E2") "  -0.05%  "
Set-TextValue $ws.Range("D3") "1.567.46"
Set-TextValue $ws.Range("E3") "  +0.15%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.12%  "
Set-TextValue $ws.Range("D5") "1.002"
Set-TextValue $ws.Range("E5") "  -0.04%  "
Set-TextValue $ws.Range("D6") "291.18"
Set-TextValue $ws.Range("E6") "  +0.22%  "
Set-TextValue $ws.Range("D7") "0.3768"
Set-TextValue $ws.Range("D8") "49.09"
Set-TextValue $ws.Range("E8") "  -0.42%  "
Set-TextValue $ws.Range("D9") "0.3395"
Set-TextValue $ws.Range("E9") "  +0.11%  "
Set-TextValue $ws.Range("D10") "0.07560"
Set-TextValue $ws.Range("E10") "  -1.17%  "
Set-TextValue $ws.Range("E11") "  -2.76%  "
Set-TextValue $ws.Range("E12") "  +0.14%  "
Set-TextValue $ws.Range("D13") "20.98"
Set-TextValue $ws.Range("E13") "  -1.81%  "
Set-TextValue $ws.Range("D14") "5.968"
Set-TextValue $ws.Range("E14") "  -1.42%  "
Set-TextValue $ws.Range("D15") "6.916"
Set-TextValue $ws.Range("E15") "  +0.03%  "
Set-TextValue $ws.Range("D16") "1.566.29"
Set-TextValue $ws.Range("E16") "  +0.14%  "
Set-TextValue $ws.Range("D17") "0.00001126"
Set-TextValue $ws.Range("E17") "  -0.07%  "
Set-TextValue $ws.Range("D18") "89.79"
Set-TextValue $ws.Range("E18") "  -0.23%  "
Set-TextValue $ws.Range("E19") "  +0.30%  "
Set-TextValue $ws.Range("E20") "  +0.08%  "
Set-TextValue $ws.Range("D21") "16.56"
Set-TextValue $ws.Range("E21") "  +0.19%  "
Set-TextValue $ws.Range("D22") "6.195"
Set-TextValue $ws.Range("E22") "  -0.88%  "
Set-TextValue $ws.Range("D23") "11.94"
Set-TextValue $ws.Range("E23") "  -0.64%  "
Set-TextValue $ws.Range("D24") "22.377.24"
Set-TextValue $ws.Range("E24") "  +0.00%  "
Set-TextValue $ws.Range("D25") "2.379"
Set-TextValue $ws.Range("E25") "  +0.59%  "
Set-TextValue $ws.Range("D26") "2.707"
Set-TextValue $ws.Range("E26") "  -3.46%  "
Set-TextValue $ws.Range("E27") "  +0.39%  "
Set-TextValue $ws.Range("D28") "148.06"
Set-TextValue $ws.Range("E28") "  +0.93%  "
Set-TextValue $ws.Range("D29") "5.025"
Set-TextValue $ws.Range("E29") "  +0.97%  "
Set-TextValue $ws.Range("D30") "125.56"
Set-TextValue $ws.Range("E30") "  -0.01%  "
Set-TextValue $ws.Range("D31") "1.737.78"
Set-TextValue $ws.Range("E31") "  +0.09%  "
Set-TextValue $ws.Range("D32") "2.020"
Set-TextValue $ws.Range("E32") "  +0.19%  "
Set-TextValue $ws.Range("D33") "6.044"
Set-TextValue $ws.Range("E33") "  -2.52%  "
Set-TextValue $ws.Range("D34") "0.9887"
Set-TextValue $ws.Range("E34") "  -2.99%  "
Set-TextValue $ws.Range("D35") "10.06"
Set-TextValue $ws.Range("E35") "  +0.06%  "
Set-TextValue $ws.Range("E36") "  +11.15%  "
Set-TextValue $ws.Range("D37") "0.08442"
Set-TextValue $ws.Range("E37") "  -0.92%  "
Set-TextValue $ws.Range("D38") "0.02483"
Set-TextValue $ws.Range("E38") "  -2.25%  "
Set-TextValue $ws.Range("D39") "0.2289"
Set-TextValue $ws.Range("E39") "  -1.38%  "
Set-TextValue $ws.Range("D40") "0.06459"
Set-TextValue $ws.Range("E40") "  +0.04%  "
Set-TextValue $ws.Range("D41") "5.408"
Set-TextValue $ws.Range("E41") "  -1.94%  "
Set-TextValue $ws.Range("D42") "0.6300"
Set-TextValue $ws.Range("E42") "  -0.66%  "
Set-TextValue $ws.Range("D43") "11.25"
Set-TextValue $ws.Range("E43") "  -3.73%  "
Set-TextValue $ws.Range("D44") "1.002"
Set-TextValue $ws.Range("E44") "  +0.03%  "
Set-TextValue $ws.Range("E45") "  -0.68%  "
Set-TextValue $ws.Range("D46") "3.799"
Set-TextValue $ws.Range("D47") "0.5921"
Set-TextValue $ws.Range("E47") "  -0.92%  "
Set-TextValue $ws.Range("D48") "2.071"
Set-TextValue $ws.Range("E48") "  -1.47%  "
Set-TextValue $ws.Range("D49") "1.261"
Set-TextValue $ws.Range("E49") "  -0.47%  "
Set-TextValue $ws.Range("D50") "124.65"
Set-TextValue $ws.Range("E50") "  -0.06%  "
Set-TextValue $ws.Range("D51") "0.07349"
Set-TextValue $ws.Range("E51") "  +1.12%  "
